# Error fix for the IFRS company list: re-key the historical-year rows (2019
# and earlier annual columns shifted) with the corrected figures, drop the
# now-unused "당기순이익(비지배)" (J) and "자본총계(비지배)" (O) columns for the
# rows that still had them, and blank out the estimate rows (2019E/2020E/2021E)
# that no longer have reported financials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected figures per row (column letter -> new value). Columns J and O are
# intentionally absent below and get their contents cleared afterwards.
$rowData = @{
    2 = @{
        D=1371; E=34; F=34; G=6; H=3; I=3; K=1055; L=679; M=376; N=376;
        P=187; Q=31; R=-15; S=-35; T=66; U=-35; V=559;
        W=2.51; X=0.22; Y=0.85; Z=0.29; AA=180.59; AB=96.2; AC=8; AD=143.51;
        AE=1006; AF=1.15; AG=0; AH=0; AI=0; AJ=37372340
    }
    3 = @{
        D=1279; E=42; F=42; G=18; H=16; I=16; K=970; L=582; M=388; N=388;
        P=187; Q=149; R=-29; S=-105; T=16; U=133; V=450;
        W=3.3; X=1.26; Y=4.23; Z=1.6; AA=150.02; AB=102.59; AC=43; AD=36.74;
        AE=1039; AF=1.53; AG=0; AH=0; AI=0; AJ=37372340
    }
    4 = @{
        D=979; E=-48; F=-48; G=-59; H=-51; I=-51; K=1021; L=432; M=589; N=589;
        P=231; Q=31; R=-156; S=194; T=23; U=8; V=312;
        W=-4.89; X=-5.22; Y=-10.46; Z=-5.14; AA=73.29; AB=151.67; AC=-128; AD=-26.54;
        AE=1277; AF=2.67; AG=0; AH=0; AI=0; AJ=46144269
    }
    5 = @{
        D=1121; E=-32; F=-32; G=44; H=41; I=41; K=2021; L=758; M=1263; N=1263;
        P=262; Q=-39; R=-256; S=397; T=22; U=-60; V=504;
        W=-2.89; X=3.67; Y=4.44; Z=2.7; AA=60.04; AB=226.23; AC=82; AD=56.25;
        AE=2413; AF=1.91; AG=0; AH=0; AI=0; AJ=52332155
    }
    6 = @{
        D=911; E=-184; F=-184; G=-316; H=-300; I=-300; K=1942; L=867; M=1074; N=1074;
        P=283; Q=-21; R=-299; S=423; T=32; U=-53; V=668;
        W=-20.21; X=-32.91; Y=-25.65; Z=-15.13; AA=80.72; AB=179.58; AC=-538; AD=-13.52;
        AE=1900; AF=3.83; AG=0; AH=0; AI=0; AJ=56551602
    }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}

# Rows 2-5 still carried the retired J (당기순이익(비지배)) and O (자본총계(비지배))
# columns from before the correction; remove them so the row matches the new
# (narrower) layout used by rows 6-9.
foreach ($r in 2..5) {
    $ws.Range("J$r").ClearContents()
    $ws.Range("O$r").ClearContents()
}

# Rows 7-9 (2020E/2021E estimate rows) no longer report any financials beyond
# the identifying columns A:C.
foreach ($r in 7..9) {
    $ws.Range("D$r`:AJ$r").ClearContents()
}
